# GBDS JANUARY FILES 2026 | fliqlo@GBDS
# Price update: FLAVORED BEER 330ml variants (rows 17-21, column C)
# SELLING PRICE per CASE goes from 205 to 229. The dependent columns
# E (CAPITAL/BOTTLE = C/D) and K (HALF CASE = C/2) are formulas and
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C17:C21").Value = 229

# Leave the cursor where the author left it when they saved the file.
$ws.Activate()
$ws.Range("C20").Select()
